$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("MATEMATICA DISCRETA", 44),
    @("METODOS MATEMATICOS I", 130),
    @("PROGRAMACION", 121),
    @("SEMINARIO DE SOLUCION DE PROBLEMAS DE PROGRAMACION", 129),
    @("SEMINARIO DE SOLUCION DE PROBLEMAS DE METODOS MATEMATICOS I", 166),
    @("ALGORITMIA", 57),
    @("ESTRUCTURAS DE DATOS I", 1),
    @("METODOS MATEMATICOS II", 109),
    @("SEMINARIO DE SOLUCION DE PROBLEMAS DE ALGORITMIA", 106),
    @("SEMINARIO DE SOLUCION DE PROBLEMAS DE ESTRUCTURAS DE DATOS I", 61),
    @("SEMINARIO DE SOLUCION DE PROBLEMAS DE METODOS MATEMATICOS II", 129),
    @("ADMINISTRACION DE REDES", 74),
    @("ESTADISTICA Y PROCESOS ESTOCASTICOS", 66),
    @("ESTRUCTURAS DE DATOS II", 11),
    @("SEMINARIO DE SOLUCION DE PROBLEMAS DE ESTRUCTURAS DE DATOS II", 80),
    @("TEORIA DE LA COMPUTACION", 79),
    @("ADMINISTRACION DE SERVIDORES", 87),
    @("BASES DE DATOS", 63),
    @("HIPERMEDIA", 69),
    @("INGENIERIA DE SOFTWARE I", 30),
    @("SEMINARIO DE SOLUCION DE PROBLEMAS DE BASES DE DATOS", 116),
    @("CONTROL DE PROYECTOS", 88),
    @("INGENIERIA DE SOFTWARE II", 82),
    @("PROGRAMACION PARA INTERNET", 169),
    @("SEGURIDAD DE LA INFORMACION", 119),
    @("SEMINARIO DE SOLUCION DE PROBLEMAS DE INGENIERIA DE SOFTWARE I", 140),
    @("SEMINARIO DE SOLUCION DE PROBLEMAS DE USO, ADAPTACION, EXPLOTACION DE SISTEMAS OPERATIVOS", 122),
    @("USO, ADAPTACION Y EXPLOTACION DE SISTEMAS OPERATIVOS", 119),
    @("ADMINISTRACION DE BASES DE DATOS", 101),
    @("ALMACENES DE DATOS (DATA WAREHOUSE)", 75),
    @("MINERIA DE DATOS", 101),
    @("CLASIFICACION INTELIGENTE DE DATOS", 88),
    @("SEMINARIO DE SOLUCION DE PROBLEMAS DE SISTEMAS BASADOS EN CONOCIMIENTO", 109),
    @("SISTEMAS BASADOS EN CONOCIMIENTO", 113)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
